$d = $word.ActiveDocument

# 1. Merge the split "Recuperació de contrasenya i Autenticació Social" runs in the TOC
#    into a single run of plain text (simulating a TOC refresh).
$d.Content.Find.Execute("Recuperació de contrasenya i Auten", $true, $false, $false, $false, $false, $true, 1, $false, "Recuperació de contrasenya i Autenticació Social", 2)
